# Insert a new record row at row 293 (shifting the existing rows 293:337
# down to 294:338) and populate the new row with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 293..337 down by one row.
$ws.Rows.Item(293).Insert()

# Populate the newly inserted row 293 with the new record's data.
$ws.Cells.Item(293, 1).Value  = 10
$ws.Cells.Item(293, 2).Value  = 'Vega Modelo de Temuco'
$ws.Cells.Item(293, 3).Value  = 'La Araucanía'
$ws.Cells.Item(293, 4).Value  = 45154
$ws.Cells.Item(293, 5).Value  = 9
$ws.Cells.Item(293, 6).Value  = 100112013
$ws.Cells.Item(293, 7).Value  = 'Alcachofa'
$ws.Cells.Item(293, 8).Value  = 'Madrigal'
$ws.Cells.Item(293, 9).Value  = 'Extra'
$ws.Cells.Item(293, 10).Value = 80
$ws.Cells.Item(293, 11).Value = 13000
$ws.Cells.Item(293, 12).Value = 13000
$ws.Cells.Item(293, 13).Value = 13000
$ws.Cells.Item(293, 14).Value = '$/caja 35 unidades'
$ws.Cells.Item(293, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(293, 16).Value = 371
$ws.Cells.Item(293, 17).Value = 35
$ws.Cells.Item(293, 18).Value = 'Hortaliza'
